$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 updates
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "643,531,686,575"

# Row 17 updates
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "794,481,830,526"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "0.72"
